# Inventory (estoque) value corrections — recompute qty/value pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: unit-rounding adjustment to the extended value only
$ws.Range("H20").Value = 98903.18

# Row 24: quantity (F) and extended value (H) recomputed
$ws.Range("F24").Value = 1059.000
$ws.Range("H24").Value = 9631.59

# Row 26: quantity (F) and extended value (H) recomputed
$ws.Range("F26").Value = 1244.000
$ws.Range("H26").Value = 12790.50

# Row 39: unit-rounding adjustment to the extended value only
$ws.Range("H39").Value = 49308.34

# Row 75: quantity (F) and extended value (H) recomputed
$ws.Range("F75").Value = 533.000
$ws.Range("H75").Value = 9980.43

# Row 76: quantity (F) and extended value (H) recomputed
$ws.Range("F76").Value = 20604.000
$ws.Range("H76").Value = 31142.16

# Row 88: quantity (F) and extended value (H) recomputed
$ws.Range("F88").Value = 237.000
$ws.Range("H88").Value = 4801.62

# Row 102: quantity (F) and extended value (H) recomputed
$ws.Range("F102").Value = 10438.800
$ws.Range("H102").Value = 23692.57

# Row 111: quantity (F) and extended value (H) recomputed
$ws.Range("F111").Value = 2187.000
$ws.Range("H111").Value = 5490.15
